{"js": "// Remove the two floating chart/scatter-plot pictures that were anchored\n// inside the document's first (and only) paragraph. Office.js treats\n// these wp:anchor drawings as floating \"shapes\" (they are not inline\n// pictures), so they are reached through body.shapes rather than\n// body.inlinePictures.\nconst shapes = context.document.body.shapes;\nshapes.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < shapes.items.length; i++) {\n  shapes.items[i].delete();\n}\nawait context.sync();\n", "ps1": "# Remove the two floating chart/scatter-plot pictures that were anchored\n# inside the document's first (and only) paragraph. These are floating\n# (anchored) drawings, so they live on $d.Shapes rather than\n# $d.InlineShapes; delete them from the end so indices stay valid.\n$d = $word.ActiveDocument\n\nfor ($i = $d.Shapes.Count; $i -ge 1; $i--) {\n    $d.Shapes.Item($i).Delete()\n}\n"}
